# Insert two new rows at 197-198 (pushing the existing rows 197:224 down to
# 199:226) and populate the two new rows with the new weekly price records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A197:A198").EntireRow.Insert()

# New row 197
$ws.Range("A197").Value = 9
$ws.Range("B197").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C197").Value = "Metropolitana"
$ws.Range("D197").Value = 44504
$ws.Range("E197").Value = 13
$ws.Range("F197").Value = 100112032
$ws.Range("G197").Value = "Zapallo italiano"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 210
$ws.Range("K197").Value = 10000
$ws.Range("L197").Value = 11000
$ws.Range("M197").Value = 10500
$ws.Range("N197").Value = "$/caja 50 unidades"
$ws.Range("O197").Value = "Región de Arica y Parinacota"
$ws.Range("P197").Value = 210
$ws.Range("Q197").Value = 50
$ws.Range("R197").Value = "Hortaliza"

# New row 198
$ws.Range("A198").Value = 9
$ws.Range("B198").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C198").Value = "Metropolitana"
$ws.Range("D198").Value = 44504
$ws.Range("E198").Value = 13
$ws.Range("F198").Value = 100112032
$ws.Range("G198").Value = "Zapallo italiano"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Segunda"
$ws.Range("J198").Value = 79
$ws.Range("K198").Value = 8000
$ws.Range("L198").Value = 8000
$ws.Range("M198").Value = 8000
$ws.Range("N198").Value = "$/caja 100 unidades"
$ws.Range("O198").Value = "Región de Arica y Parinacota"
$ws.Range("P198").Value = 80
$ws.Range("Q198").Value = 100
$ws.Range("R198").Value = "Hortaliza"
